$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches original inlineStr cell type) for numeric-looking
# Price values, so Excel does not auto-convert them to numbers on assignment.
$textCells = @(
    'D5',
    'D6',
    'D9',
    'D10',
    'D11',
    'D12',
    'D13',
    'D14',
    'D17',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D24',
    'D25',
    'D26',
    'D27',
    'D28',
    'D29',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D50',
    'D51'
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '64.311.76'
$ws.Range('E2').Value = '  +0.39%  '

$ws.Range('D3').Value = '3.129.35'
$ws.Range('E3').Value = '  +2.89%  '

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = '560.31'
$ws.Range('E5').Value = '  +2.23%  '

$ws.Range('D6').Value = '144.07'
$ws.Range('E6').Value = '  +4.30%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '3.123.17'
$ws.Range('E8').Value = '  +3.08%  '

$ws.Range('D9').Value = '0.493'
$ws.Range('E9').Value = '  +1.72%  '

$ws.Range('D10').Value = '6.70'
$ws.Range('E10').Value = '  +3.92%  '

$ws.Range('D11').Value = '0.153'
$ws.Range('E11').Value = '  +0.76%  '

$ws.Range('D12').Value = '0.465'
$ws.Range('E12').Value = '  +2.19%  '

$ws.Range('D13').Value = '36.53'
$ws.Range('E13').Value = '  +2.66%  '

$ws.Range('D14').Value = '0.0000221'
$ws.Range('E14').Value = '  +1.58%  '

$ws.Range('D15').Value = '3.635.51'
$ws.Range('E15').Value = '  +3.06%  '

$ws.Range('D16').Value = '64.364.11'
$ws.Range('E16').Value = '  +0.50%  '

$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '0.112'
$ws.Range('E17').Value = '  +1.18%  '

$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.127.03'
$ws.Range('E18').Value = '  +2.90%  '

$ws.Range('D19').Value = '514.11'
$ws.Range('E19').Value = '  +6.03%  '

$ws.Range('D20').Value = '6.82'
$ws.Range('E20').Value = '  +4.12%  '

$ws.Range('D21').Value = '13.96'
$ws.Range('E21').Value = '  +2.84%  '

$ws.Range('D22').Value = '0.713'
$ws.Range('E22').Value = '  +4.84%  '

$ws.Range('D23').Value = '7.43'
$ws.Range('E23').Value = '  +4.64%  '

$ws.Range('D24').Value = '12.89'
$ws.Range('E24').Value = '  +4.54%  '

$ws.Range('D25').Value = '78.42'
$ws.Range('E25').Value = '  +0.76%  '

$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.11%  '

$ws.Range('D27').Value = '8.92'
$ws.Range('E27').Value = '  +16.34%  '

$ws.Range('D28').Value = '2.82'
$ws.Range('E28').Value = '  +4.69%  '

$ws.Range('D29').Value = '2.15'
$ws.Range('E29').Value = '  +4.71%  '

$ws.Range('E30').Value = '  +0.07%  '

$ws.Range('D31').Value = '26.40'
$ws.Range('E31').Value = '  +2.63%  '

$ws.Range('D32').Value = '2.59'
$ws.Range('E32').Value = '  -1.51%  '

$ws.Range('D33').Value = '1.13'
$ws.Range('E33').Value = '  +2.75%  '

$ws.Range('D34').Value = '544.35'
$ws.Range('E34').Value = '  -7.42%  '

$ws.Range('D35').Value = '5.37'
$ws.Range('E35').Value = '  -0.05%  '

$ws.Range('D36').Value = '6.05'
$ws.Range('E36').Value = '  +3.61%  '

$ws.Range('D37').Value = '53.87'
$ws.Range('E37').Value = '  +4.10%  '

$ws.Range('D38').Value = '0.0432'
$ws.Range('E38').Value = '  +6.22%  '

$ws.Range('D39').Value = '0.0823'
$ws.Range('E39').Value = '  +4.42%  '

$ws.Range('D40').Value = '3.129.38'
$ws.Range('E40').Value = '  +7.30%  '

$ws.Range('D41').Value = '0.122'
$ws.Range('E41').Value = '  +2.83%  '

$ws.Range('D42').Value = '2.75'
$ws.Range('E42').Value = '  -1.67%  '

$ws.Range('D43').Value = '8.24'
$ws.Range('E43').Value = '  +0.86%  '

$ws.Range('D44').Value = '0.268'
$ws.Range('E44').Value = '  +11.47%  '

$ws.Range('D45').Value = '2.21'
$ws.Range('E45').Value = '  +7.03%  '

$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.03%  '

$ws.Range('D47').Value = '25.42'
$ws.Range('E47').Value = '  +3.42%  '

$ws.Range('D48').Value = '121.01'
$ws.Range('E48').Value = '  +2.12%  '

$ws.Range('D49').Value = '0.0₃0520'
$ws.Range('E49').Value = '  -1.78%  '

$ws.Range('D50').Value = '0.108'
$ws.Range('E50').Value = '  +0.33%  '

$ws.Range('D51').Value = '2.10'
$ws.Range('E51').Value = '  +3.61%  '
